# Delete rows corresponding to even_MAG-GUT28366.fa, even_MAG-GUT83851.fa,
# and even_MAG-GUT85090.fa (the "(reject)" rows), leaving the remaining
# rows shifted up. Original rows 5 and 6 are deleted first (from the
# bottom up) to keep row indices stable, then row 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(3).Delete()
